# Weekly CompStat (CS-EN-US-PBQS) update: new crime data collected.
# Updates the report header (issue number + reporting week dates) and
# refreshes all the weekly/28-day/YTD/2-year/15-year/32-year crime figures
# for rows 14-33 of the CompStat worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header text -------------------------------------------------
# A8:  "Volume 32   Number  48"  -> "Volume 32   Number  49"
# C9:  "Report Covering the Week  11/24/2025  Through  11/30/2025"
#      -> "Report Covering the Week  12/1/2025  Through  12/7/2025"
$ws.Range("A8").Value2 = "Volume 32   Number  49"
$ws.Range("C9").Value2 = "Report Covering the Week  12/1/2025  Through  12/7/2025"

# --- Weekly crime statistics table (rows 14-33) --------------------------
$ws.Range("J14").Value2 = 36
$ws.Range("K14").Value2 = -16.666666666666
$ws.Range("L14").Value2 = 30.434782608695
$ws.Range("M14").Value2 = -54.545454545454
$ws.Range("N14").Value2 = -79.591836734693
$ws.Range("D15").Value2 = 4
$ws.Range("E15").Value2 = -75
$ws.Range("F15").Value2 = 11
$ws.Range("G15").Value2 = 14
$ws.Range("H15").Value2 = -21.428571428571
$ws.Range("J15").Value2 = 186
$ws.Range("K15").Value2 = 13.978494623655
$ws.Range("L15").Value2 = 37.662337662337
$ws.Range("M15").Value2 = 73.770491803278
$ws.Range("N15").Value2 = -35.562310030395
$ws.Range("C16").Value2 = 20
$ws.Range("D16").Value2 = 23
$ws.Range("E16").Value2 = -13.043478260869
$ws.Range("F16").Value2 = 93
$ws.Range("G16").Value2 = 80
$ws.Range("H16").Value2 = 16.25
$ws.Range("I16").Value2 = 1171
$ws.Range("J16").Value2 = 1310
$ws.Range("K16").Value2 = -10.610687022900
$ws.Range("L16").Value2 = -8.372456964006
$ws.Range("M16").Value2 = -41.799204771371
$ws.Range("N16").Value2 = -84.338638491373
$ws.Range("C17").Value2 = 47
$ws.Range("D17").Value2 = 51
$ws.Range("E17").Value2 = -7.843137254901
$ws.Range("F17").Value2 = 200
$ws.Range("G17").Value2 = 209
$ws.Range("H17").Value2 = -4.306220095693
$ws.Range("I17").Value2 = 2846
$ws.Range("J17").Value2 = 2954
$ws.Range("K17").Value2 = -3.656059580230
$ws.Range("L17").Value2 = 5.996275605214
$ws.Range("M17").Value2 = 76.222910216718
$ws.Range("N17").Value2 = -20.369334079462
$ws.Range("C18").Value2 = 30
$ws.Range("D18").Value2 = 27
$ws.Range("E18").Value2 = 11.111111111111
$ws.Range("F18").Value2 = 83
$ws.Range("G18").Value2 = 99
$ws.Range("H18").Value2 = -16.161616161616
$ws.Range("I18").Value2 = 993
$ws.Range("J18").Value2 = 1032
$ws.Range("K18").Value2 = -3.779069767441
$ws.Range("L18").Value2 = -7.196261682242
$ws.Range("M18").Value2 = -51.889534883720
$ws.Range("N18").Value2 = -88.558589699274
$ws.Range("C19").Value2 = 77
$ws.Range("D19").Value2 = 56
$ws.Range("E19").Value2 = 37.5
$ws.Range("F19").Value2 = 273
$ws.Range("G19").Value2 = 239
$ws.Range("H19").Value2 = 14.225941422594
$ws.Range("I19").Value2 = 3292
$ws.Range("J19").Value2 = 3105
$ws.Range("K19").Value2 = 6.022544283413
$ws.Range("L19").Value2 = -3.714536414156
$ws.Range("M19").Value2 = 19.361856417694
$ws.Range("N19").Value2 = -56.478053939714
$ws.Range("C20").Value2 = 19
$ws.Range("D20").Value2 = 35
$ws.Range("E20").Value2 = -45.714285714285
$ws.Range("F20").Value2 = 107
$ws.Range("G20").Value2 = 153
$ws.Range("H20").Value2 = -30.065359477124
$ws.Range("I20").Value2 = 1784
$ws.Range("J20").Value2 = 1887
$ws.Range("K20").Value2 = -5.458399576046
$ws.Range("L20").Value2 = 3.420289855072
$ws.Range("M20").Value2 = 7.859733978234
$ws.Range("N20").Value2 = -90.202108963093
$ws.Range("C21").Value2 = 194
$ws.Range("D21").Value2 = 197
$ws.Range("E21").Value2 = -1.522842639593
$ws.Range("F21").Value2 = 767
$ws.Range("G21").Value2 = 796
$ws.Range("H21").Value2 = -3.643216080402
$ws.Range("I21").Value2 = 10328
$ws.Range("J21").Value2 = 10510
$ws.Range("K21").Value2 = -1.731684110371
$ws.Range("L21").Value2 = -0.251110681862
$ws.Range("M21").Value2 = 0.359537459916
$ws.Range("N21").Value2 = -77.537082952716
$ws.Range("C22").Value2 = 3
$ws.Range("E22").Value2 = 50
$ws.Range("F22").Value2 = 6
$ws.Range("H22").Value2 = -33.333333333333
$ws.Range("I22").Value2 = 109
$ws.Range("J22").Value2 = 126
$ws.Range("K22").Value2 = -13.492063492063
$ws.Range("L22").Value2 = -7.627118644067
$ws.Range("M22").Value2 = -2.678571428571
$ws.Range("D23").Value2 = 2
$ws.Range("E23").Value2 = 100
$ws.Range("F23").Value2 = 16
$ws.Range("H23").Value2 = 100
$ws.Range("I23").Value2 = 214
$ws.Range("J23").Value2 = 200
$ws.Range("K23").Value2 = 7
$ws.Range("L23").Value2 = -3.603603603603
$ws.Range("M23").Value2 = 50.704225352112
$ws.Range("C24").Value2 = 186
$ws.Range("D24").Value2 = 205
$ws.Range("E24").Value2 = -9.268292682926
$ws.Range("F24").Value2 = 800
$ws.Range("G24").Value2 = 777
$ws.Range("H24").Value2 = 2.960102960102
$ws.Range("I24").Value2 = 8769
$ws.Range("J24").Value2 = 8908
$ws.Range("K24").Value2 = -1.560395150426
$ws.Range("L24").Value2 = 1.669565217391
$ws.Range("M24").Value2 = 51.476939022283
$ws.Range("C25").Value2 = 68
$ws.Range("D25").Value2 = 83
$ws.Range("E25").Value2 = -18.072289156626
$ws.Range("F25").Value2 = 349
$ws.Range("G25").Value2 = 374
$ws.Range("H25").Value2 = -6.684491978609
$ws.Range("I25").Value2 = 3514
$ws.Range("J25").Value2 = 3956
$ws.Range("K25").Value2 = -11.172901921132
$ws.Range("L25").Value2 = 17.643120187479
$ws.Range("C26").Value2 = 87
$ws.Range("D26").Value2 = 108
$ws.Range("E26").Value2 = -19.444444444444
$ws.Range("F26").Value2 = 348
$ws.Range("G26").Value2 = 400
$ws.Range("H26").Value2 = -13
$ws.Range("I26").Value2 = 4612
$ws.Range("J26").Value2 = 4903
$ws.Range("K26").Value2 = -5.935141749949
$ws.Range("L26").Value2 = 7.205950720595
$ws.Range("M26").Value2 = 9.861838970938
$ws.Range("D27").Value2 = 6
$ws.Range("E27").Value2 = -66.666666666666
$ws.Range("F27").Value2 = 13
$ws.Range("G27").Value2 = 22
$ws.Range("H27").Value2 = -40.909090909090
$ws.Range("I27").Value2 = 264
$ws.Range("J27").Value2 = 285
$ws.Range("K27").Value2 = -7.368421052631
$ws.Range("L27").Value2 = 1.930501930501
$ws.Range("C28").Value2 = 6
$ws.Range("D28").Value2 = 7
$ws.Range("E28").Value2 = -14.285714285714
$ws.Range("F28").Value2 = 40
$ws.Range("G28").Value2 = 35
$ws.Range("H28").Value2 = 14.285714285714
$ws.Range("I28").Value2 = 431
$ws.Range("J28").Value2 = 415
$ws.Range("K28").Value2 = 3.855421686746
$ws.Range("L28").Value2 = 19.060773480663
$ws.Range("D29").Value2 = 3
$ws.Range("G29").Value2 = 8
$ws.Range("H29").Value2 = -50
$ws.Range("J29").Value2 = 110
$ws.Range("K29").Value2 = -40.909090909090
$ws.Range("L29").Value2 = -23.529411764705
$ws.Range("M29").Value2 = -65.789473684210
$ws.Range("N29").Value2 = -87.077534791252
$ws.Range("D30").Value2 = 3
$ws.Range("G30").Value2 = 7
$ws.Range("H30").Value2 = -57.142857142857
$ws.Range("J30").Value2 = 90
$ws.Range("K30").Value2 = -43.333333333333
$ws.Range("L30").Value2 = -22.727272727272
$ws.Range("M30").Value2 = -65.771812080536
$ws.Range("N30").Value2 = -88.864628820960
$ws.Range("L31").Value2 = -34.545454545454
$ws.Range("G33").Value2 = 1
$ws.Range("H33").Value2 = 100
$ws.Range("I33").Value2 = 30
$ws.Range("K33").Value2 = 30.434782608695
$ws.Range("L33").Value2 = -21.052631578947

# --- Cells that change between numeric and "N/A" text representation ----
# Hate Crimes (row 31): 28-day 2025 count and %chg become not meaningful
# ("0" / "***.*" placeholder text), replacing the prior numeric 1 / 0.
$ws.Range("D31").Copy()
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value2 = "0"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value2 = "***.*"
$ws.Range("G31:H31").PasteSpecial(-4122)

# Traffic Fatalities (row 33): Week-to-Date 2025 and 28-Day 2025 columns
# now have reported numeric counts instead of the "N/A" placeholder text.
$ws.Range("G33").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("C33").Value2 = 1
$ws.Range("F33").PasteSpecial(-4122)
$ws.Range("F33").Value2 = 2
$excel.CutCopyMode = $false
